$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: Table-S3a-GO-pval-0.05-DE -> Table-S5a-GO-pval-0.05-DE
$ws.Name = "Table-S5a-GO-pval-0.05-DE"

# Bold the header row (A1:L1) - introduces a new bold font/style used only by the header
$ws.Range("A1:L1").Font.Bold = $true

# Move/restore the active selection to E31
[void]$ws.Range("E31").Select()
